# Update "ad_fcv" worksheet: add std-dev columns for each horizon, rename
# CART -> DTREE, and drop the NB algorithm row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the NB row (was row 9, algorithm index 7) entirely, so SVM's
#    data collapses up into row 8.
# ---------------------------------------------------------------------
$ws.Rows("9:9").Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. Header row - B1 ("Algorithm") is unchanged. Existing mean headers
#    (C1:G1) get a " mean" suffix, each followed by a new " std" column,
#    extending the header row out to column L.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# Copy the bold/border/centered header style from an existing header cell
# (B1) onto the newly added header cells H1:L1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:L1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fix up algorithm labels after the row shift:
#      - row 5 (index 3): CART -> DTREE
#      - row 8 (index 6): was "NB" before the delete, now needs to read
#        "SVM" (the data that shifted up from the deleted row's sibling)
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "DTREE"
$ws.Range("B8").Value = "SVM"

# ---------------------------------------------------------------------
# 4. Rewrite the numeric data block (rows 2-8, columns C:L) with the new
#    mean/std values. Column B (algorithm names, other than the DTREE
#    rename above) and column A (0-based index) are unchanged.
# ---------------------------------------------------------------------

# Row 2 - LR
$ws.Range("C2").Value = 0.8264268906274909
$ws.Range("D2").Value = 0.01122763501421739
$ws.Range("E2").Value = 0.8143328379751006
$ws.Range("F2").Value = 0.02012404106331644
$ws.Range("G2").Value = 0.8109150785487065
$ws.Range("H2").Value = 0.02577186844924204
$ws.Range("I2").Value = 0.8012891376516655
$ws.Range("J2").Value = 0.02539177282170518
$ws.Range("K2").Value = 0.7957997515833358
$ws.Range("L2").Value = 0.01661393786092331

# Row 3 - LDA
$ws.Range("C3").Value = 0.828606625945117
$ws.Range("D3").Value = 0.01555913010989755
$ws.Range("E3").Value = 0.8157343673018017
$ws.Range("F3").Value = 0.018806521292397
$ws.Range("G3").Value = 0.8085701086391921
$ws.Range("H3").Value = 0.02186849292478791
$ws.Range("I3").Value = 0.7975634558053991
$ws.Range("J3").Value = 0.02544564488415158
$ws.Range("K3").Value = 0.7971571301335163
$ws.Range("L3").Value = 0.01789072582137708

# Row 4 - KNN
$ws.Range("C4").Value = 0.7849732995398155
$ws.Range("D4").Value = 0.01402292014079909
$ws.Range("E4").Value = 0.7860542918352966
$ws.Range("F4").Value = 0.01293600771969228
$ws.Range("G4").Value = 0.7913651939779491
$ws.Range("H4").Value = 0.0220138601406211
$ws.Range("I4").Value = 0.7859606125915753
$ws.Range("J4").Value = 0.02393877279497811
$ws.Range("K4").Value = 0.7899518282178816
$ws.Range("L4").Value = 0.01780104247844531

# Row 5 - DTREE
$ws.Range("C5").Value = 0.7785419873187778
$ws.Range("D5").Value = 0.01996766690408337
$ws.Range("E5").Value = 0.7583743501208782
$ws.Range("F5").Value = 0.02595044310704902
$ws.Range("G5").Value = 0.7627610847515232
$ws.Range("H5").Value = 0.01893351355317181
$ws.Range("I5").Value = 0.7572816705829635
$ws.Range("J5").Value = 0.02150987817729219
$ws.Range("K5").Value = 0.7787842209150448
$ws.Range("L5").Value = 0.047325556851765

# Row 6 - RTREE
$ws.Range("C6").Value = 0.8099408374985577
$ws.Range("D6").Value = 0.01450498083821188
$ws.Range("E6").Value = 0.7868341496122632
$ws.Range("F6").Value = 0.01774035922152552
$ws.Range("G6").Value = 0.7658460597753765
$ws.Range("H6").Value = 0.01859032752324731
$ws.Range("I6").Value = 0.7498881417320469
$ws.Range("J6").Value = 0.02538888135000571
$ws.Range("K6").Value = 0.7272122875408982
$ws.Range("L6").Value = 0.02487460848362794

# Row 7 - XTREE
$ws.Range("C7").Value = 0.8385617961403982
$ws.Range("D7").Value = 0.01595641532629525
$ws.Range("E7").Value = 0.8146635517660957
$ws.Range("F7").Value = 0.02239503821140553
$ws.Range("G7").Value = 0.7996837897948321
$ws.Range("H7").Value = 0.01806656889432205
$ws.Range("I7").Value = 0.7926334592466758
$ws.Range("J7").Value = 0.03279898827721388
$ws.Range("K7").Value = 0.7921164679670845
$ws.Range("L7").Value = 0.03028093743184303

# Row 8 - SVM (index 6; NB row removed above, SVM shifted up)
$ws.Range("C8").Value = 0.8300972622867798
$ws.Range("D8").Value = 0.01371636570300042
$ws.Range("E8").Value = 0.8218032467719066
$ws.Range("F8").Value = 0.01504539490575349
$ws.Range("G8").Value = 0.8186460356550274
$ws.Range("H8").Value = 0.02423171630635305
$ws.Range("I8").Value = 0.8121240718249185
$ws.Range("J8").Value = 0.02198892462307928
$ws.Range("K8").Value = 0.8016338676450649
$ws.Range("L8").Value = 0.02048257982777268

Write-Host "Edit complete"
